$wb = $excel.ActiveWorkbook

# --- Login sheet: update selection to the whole of row 2 ---
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Rows.Item(2).Select() | Out-Null

# --- New "Parameter" sheet (extra login fixture rows) ---
$wsParam = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsParam.Name = "Parameter"

$wsParam.Columns.Item(1).ColumnWidth = 12.166666666666666
$wsParam.Columns.Item(2).ColumnWidth = 20.666666666666668

$wsParam.Range("A1").Value = "userName"
$wsParam.Range("B1").Value = "passWord"

$wsParam.Hyperlinks.Add($wsParam.Range("A2"), "mailto:bbb@bbb.com", "", "", "bbb@bbb.com") | Out-Null
$wsParam.Range("B2").Value = "bbbbb"

$wsParam.Hyperlinks.Add($wsParam.Range("A3"), "mailto:ccc@ccc.com", "", "", "ccc@ccc.com") | Out-Null
$wsParam.Range("B3").Value = "ccccc"

$wsParam.Hyperlinks.Add($wsParam.Range("A4"), "mailto:eee@eee.com", "", "", "eee@eee.com") | Out-Null
$wsParam.Range("B4").Value = "eeeee"

$wsParam.Range("A2").Style = "Hyperlink"
$wsParam.Range("A3").Style = "Hyperlink"
$wsParam.Range("A4").Style = "Hyperlink"

$wsParam.Range("B13").Select() | Out-Null

# --- New "test_suite" sheet (run-mode control table) ---
$wsSuite = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsSuite.Name = "test_suite"

$wsSuite.Columns.Item(1).ColumnWidth = 15.666666666666666
$wsSuite.Columns.Item(2).ColumnWidth = 19.998697916666668

$wsSuite.Range("A1").Value = "TCID"
$wsSuite.Range("B1").Value = "Runmode"

$wsSuite.Range("A2").Value = "Login"
$wsSuite.Range("A3").Value = "Ordering"

$wsSuite.Range("A4").Font.Name = "Helvetica"
$wsSuite.Range("A4").Value = "Parameter"

$wsSuite.Range("A5").Font.Name = "Helvetica"
$wsSuite.Range("A5").Value = "VerifyLoginPage"

$wsSuite.Range("B2").Value = "N"
$wsSuite.Range("B3").Value = "N"
$wsSuite.Range("B4").Value = "Y"
$wsSuite.Range("B5").Value = "N"

$wsSuite.Range("B12").Select() | Out-Null
